# Re-run of the TM52 test-job export: the "readme" index table is
# regenerated with a different column order (Author moved next to the
# index column, Date moved to the end) and a new run date/time stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")

# --- readme!Table1 header row: index, Author, sheet_name, JobNo, Date ---
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "sheet_name"
$ws.Range("D1").Value = "JobNo"
$ws.Range("E1").Value = "Date"

# --- readme!Table1 body rows: re-populate per the new column layout ---
$sheetNames = @(
    "Project Information",
    "Criterion % Definitions",
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

# Keep the Date column a text cell (matches the original "20220224"
# string cell) rather than letting the numeric-looking value auto-convert
# to a real number.
$ws.Range("E2:E12").NumberFormat = "@"

for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = "jovyan"
    $ws.Cells.Item($row, 3).Value = $sheetNames[$i]
    $ws.Cells.Item($row, 4).Value = "/c/e"
    $ws.Cells.Item($row, 5).Value = "20220225"
}

# --- Project Information sheet: refresh the analysis run timestamp ---
$infoWs = $wb.Worksheets.Item("Project Information")
$infoWs.Range("B11").Value = "2022-02-25 14:06:11.166802"
